$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new data rows before the current row 135. This shifts the
# existing rows 135-138 down to 138-141, matching the new dimension
# A1:R141.
$ws.Rows.Item(135).Insert()
$ws.Rows.Item(135).Insert()
$ws.Rows.Item(135).Insert()

function Set-DataRow {
    param(
        $Row, $MercadoId, $Mercado, $Region, $Fecha, $Codreg, $CategoriaId,
        $Categoria, $Variedad, $Calidad, $Volumen, $PrecioMin, $PrecioMax,
        $PrecioProm, $Unidad, $Origen, $PrecioKg, $KgUnidades, $Clasificacion
    )
    $ws.Cells.Item($Row, 1).Value = $MercadoId
    $ws.Cells.Item($Row, 2).Value = $Mercado
    $ws.Cells.Item($Row, 3).Value = $Region
    $ws.Cells.Item($Row, 4).Value = $Fecha
    $ws.Cells.Item($Row, 5).Value = $Codreg
    $ws.Cells.Item($Row, 6).Value = $CategoriaId
    $ws.Cells.Item($Row, 7).Value = $Categoria
    $ws.Cells.Item($Row, 8).Value = $Variedad
    $ws.Cells.Item($Row, 9).Value = $Calidad
    $ws.Cells.Item($Row, 10).Value = $Volumen
    $ws.Cells.Item($Row, 11).Value = $PrecioMin
    $ws.Cells.Item($Row, 12).Value = $PrecioMax
    $ws.Cells.Item($Row, 13).Value = $PrecioProm
    $ws.Cells.Item($Row, 14).Value = $Unidad
    $ws.Cells.Item($Row, 15).Value = $Origen
    $ws.Cells.Item($Row, 16).Value = $PrecioKg
    $ws.Cells.Item($Row, 17).Value = $KgUnidades
    $ws.Cells.Item($Row, 18).Value = $Clasificacion
}

Set-DataRow 135 2 "Comercializadora del Agro de Limarí" "Coquimbo" 44448 4 100112013 "Alcachofa" "Argentina(o)" "Primera" 1500 9000 10000 9500 "`$/caja 50 unidades" "Provincia de Limarí" 190 50 "Hortaliza"
Set-DataRow 136 2 "Comercializadora del Agro de Limarí" "Coquimbo" 44448 4 100112013 "Alcachofa" "Española" "Primera" 1800 10000 11000 10500 "`$/caja 30 unidades" "Provincia de Limarí" 350 30 "Hortaliza"
Set-DataRow 137 2 "Comercializadora del Agro de Limarí" "Coquimbo" 44448 4 100112013 "Alcachofa" "Española" "Segunda" 1100 8000 9000 8500 "`$/caja 40 unidades" "Provincia de Limarí" 212 40 "Hortaliza"
